$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the edited range to Text format first so Excel does not auto-convert
# numeric-looking strings (e.g. "59.926.49", "5.07") into numbers/dates.
$editRange = $ws.Range("B2:E51")
$editRange.NumberFormat = "@"

$ws.Range("D2").Value = "59.926.49"
$ws.Range("E2").Value = "  -5.14%  "
$ws.Range("D3").Value = "2.971.41"
$ws.Range("E3").Value = "  -6.84%  "
$ws.Range("D5").Value = "572.13"
$ws.Range("E5").Value = "  -3.56%  "
$ws.Range("D6").Value = "124.99"
$ws.Range("E6").Value = "  -7.82%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").Value = "2.967.84"
$ws.Range("E8").Value = "  -6.83%  "
$ws.Range("E9").Value = "  -3.16%  "
$ws.Range("E10").Value = "  -6.19%  "
$ws.Range("D11").Value = "5.07"
$ws.Range("E11").Value = "  -3.40%  "
$ws.Range("E12").Value = "  -3.91%  "
$ws.Range("D13").Value = "0.0000222"
$ws.Range("E13").Value = "  -6.51%  "
$ws.Range("D14").Value = "32.60"
$ws.Range("E14").Value = "  -6.41%  "
$ws.Range("E15").Value = "  -0.71%  "
$ws.Range("D16").Value = "3.465.01"
$ws.Range("E16").Value = "  -6.80%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "2.969.40"
$ws.Range("E17").Value = "  -6.79%  "
$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").Value = "60.002.63"
$ws.Range("E18").Value = "  -4.99%  "
$ws.Range("D19").Value = "6.18"
$ws.Range("E19").Value = "  -6.31%  "
$ws.Range("D20").Value = "432.81"
$ws.Range("E20").Value = "  -6.73%  "
$ws.Range("D21").Value = "13.05"
$ws.Range("E21").Value = "  -7.27%  "
$ws.Range("D22").Value = "0.660"
$ws.Range("E22").Value = "  -5.67%  "
$ws.Range("D23").Value = "6.99"
$ws.Range("E23").Value = "  -8.79%  "
$ws.Range("D24").Value = "12.72"
$ws.Range("E24").Value = "  -5.29%  "
$ws.Range("D25").Value = "78.85"
$ws.Range("E25").Value = "  -4.91%  "
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("D27").Value = "0.999"
$ws.Range("E27").Value = "  -0.12%  "
$ws.Range("E28").Value = "  -5.89%  "
$ws.Range("D29").Value = "7.19"
$ws.Range("E29").Value = "  -7.21%  "
$ws.Range("E30").Value = "  -7.97%  "
$ws.Range("D31").Value = "6.16"
$ws.Range("E31").Value = "  -9.37%  "
$ws.Range("D32").Value = "25.29"
$ws.Range("E32").Value = "  -7.43%  "
$ws.Range("D33").Value = "0.0929"
$ws.Range("E33").Value = "  -9.90%  "
$ws.Range("E34").Value = "  -9.27%  "
$ws.Range("D35").Value = "0.945"
$ws.Range("E35").Value = "  -8.70%  "
$ws.Range("D36").Value = "5.58"
$ws.Range("E36").Value = "  -4.55%  "
$ws.Range("D37").Value = "49.61"
$ws.Range("E37").Value = "  -3.40%  "
$ws.Range("D38").Value = "0.0₃0658"
$ws.Range("E38").Value = "  -7.49%  "
$ws.Range("E39").Value = "  -8.13%  "
$ws.Range("D40").Value = "7.89"
$ws.Range("E40").Value = "  -2.82%  "
$ws.Range("E41").Value = "  -3.52%  "
$ws.Range("D42").Value = "380.99"
$ws.Range("D43").Value = "2.47"
$ws.Range("E43").Value = "  -7.16%  "
$ws.Range("D44").Value = "2.619.73"
$ws.Range("E44").Value = "  -6.73%  "
$ws.Range("E46").Value = "  -6.99%  "
$ws.Range("E47").Value = "  -7.22%  "
$ws.Range("D48").Value = "118.39"
$ws.Range("E48").Value = "  -4.91%  "
$ws.Range("E49").Value = "  -4.83%  "
$ws.Range("D50").Value = "23.33"
$ws.Range("E50").Value = "  -7.98%  "
$ws.Range("D51").Value = "31.33"
$ws.Range("E51").Value = "  -11.32%  "

# Restore the original (default) style so only values changed, matching the source diff.
$editRange.Style = "Normal"

